$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A91").Value = "What was the outcome of your application?"
[void]$ws.Range("A91").Select()
